$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")
$ws.Range("E205").NumberFormat = "General"
